$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the header cells: "<name>_old" -> "<name>_FV2210"
#                              "<name>_new" -> "<name>_FV2304"
# ---------------------------------------------------------------------------
for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $val = $cell.Value2
    if ($val -ne $null) {
        if ($val.EndsWith("_old")) {
            $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2210"
        } elseif ($val.EndsWith("_new")) {
            $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2304"
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Turn the data range into an Excel Table ("Table1") with an autofilter.
#    The header row (A1:U1) already carries direct formatting (bold, grey
#    fill, border). Preserve that exact formatting across the ListObjects.Add
#    call by round-tripping it through an unused scratch range, since the
#    engine otherwise captures the header's look into a brand-new dxf.
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$scratchRange = $ws.Range("A1000:U1000")

$headerRange.Copy()
$scratchRange.PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$headerRange.ClearFormats()

$tableRange = $ws.Range("A1:U65")
$lo = $ws.ListObjects().Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"

$scratchRange.Copy()
$headerRange.PasteSpecial(-4122)    # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Rows("1000:1000").Delete()

# ---------------------------------------------------------------------------
# 3) Freeze the header row (split after row 1).
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
